# nested_logit_MLE_toy_model.xlsx -- "Fixing the nested logit MLE now to be
# computed with torch elements"
#
# 1) Rename Sheet1 -> regular, Sheet2 -> csv (Excel auto-rewrites the
#    Sheet1! formula references on the csv sheet to regular!).
# 2) On with_zeros: add the LN(...) helper cells (P10:R12), the new
#    row-8 SUM totals (P8:R8), and bring the L/M/N-derived LN() block
#    (P5:R7) and the SUM() rollups (G14, G19, G27, H34) to the same
#    "#,##0.00" number format.
# 3) Restore the view state: with_zeros selection moves to G27/topLeft B4
#    and is no longer the active tab; "regular" becomes the active tab.

$wb = $excel.ActiveWorkbook

$numFmt = "#,##0.00"

# ---------------------------------------------------------------------
# 1. Rename sheets (formulas on "csv" referencing Sheet1! auto-update to
#    regular! as part of the rename).
# ---------------------------------------------------------------------
$wsRegular = $wb.Worksheets.Item("Sheet1")
$wsRegular.Name = "regular"

$wsCsv = $wb.Worksheets.Item("Sheet2")
$wsCsv.Name = "csv"

$wsZeros = $wb.Worksheets.Item("with_zeros")

# ---------------------------------------------------------------------
# 2. with_zeros: new nested-logit helper formulas + number formats.
# ---------------------------------------------------------------------

# Existing LN(...)*k / LN(...) block (P5:R7) picks up the accounting
# number format used across the sheet.
$wsZeros.Range("P5:R7").NumberFormat = $numFmt

# New row 8: totals of the P6:R7 block (plus R5 for the R column).
$wsZeros.Range("P8").Formula = "=SUM(P6:P7)"
$wsZeros.Range("Q8").Formula = "=SUM(Q6:Q7)"
$wsZeros.Range("R8").Formula = "=SUM(R6:R7)+R5"
$wsZeros.Range("P8:R8").NumberFormat = $numFmt

# New rows 10-12: plain LN() of L5:N5, L6:N6, L7:N7 (the "per row" logit
# denominators, without the *3 / *4 scaling used in row 5/6).
$wsZeros.Range("P10").Formula = "=LN(L5)"
$wsZeros.Range("Q10").Formula = "=LN(M5)"
$wsZeros.Range("R10").Formula = "=LN(N5)"
$wsZeros.Range("P10:R10").NumberFormat = $numFmt

$wsZeros.Range("P11").Formula = "=LN(L6)"
$wsZeros.Range("Q11").Formula = "=LN(M6)"
$wsZeros.Range("R11").Formula = "=LN(N6)"
$wsZeros.Range("P11:R11").NumberFormat = $numFmt

$wsZeros.Range("P12").Formula = "=LN(L7)"
$wsZeros.Range("Q12").Formula = "=LN(M7)"
$wsZeros.Range("R12").Formula = "=LN(N7)"
$wsZeros.Range("P12:R12").NumberFormat = $numFmt

# The SUM() rollups keep the same "#,##0.00" look (style index shifts in
# the file because of the newly inserted cellXf, value/formula unchanged).
$wsZeros.Range("G14").NumberFormat = $numFmt
$wsZeros.Range("G19").NumberFormat = $numFmt
$wsZeros.Range("G27").NumberFormat = $numFmt
$wsZeros.Range("H34").NumberFormat = $numFmt

# ---------------------------------------------------------------------
# 3. View state: with_zeros scrolls to show row 27 and loses tab focus;
#    "regular" becomes the active sheet/tab.
# ---------------------------------------------------------------------
$wsZeros.Activate()
$wsZeros.Range("G27").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2

$wsRegular.Activate()
